$d = $word.ActiveDocument

# Add the built-in "Footnote Text" paragraph style (word/styles.xml), mirroring
# Word's behaviour when a Footnote Text style is (re)introduced into a
# document's style sheet: based on Normal, follow-on style is itself,
# ui priority 9, unhideWhenUsed + qFormat set (same shape as the
# neighbouring BlockQuote style already present in the document).
$footnoteText = $d.Styles.Add("Footnote Text", 1)
$footnoteText.BaseStyle = "Normal"
$footnoteText.NextParagraphStyle = "FootnoteText"
$footnoteText.Priority = 9
$footnoteText.UnhideWhenUsed = $true
$footnoteText.QuickStyle = $true
